$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "What's the maximum number of characters in a single text entry?"
$ws.Range("B6").Value = "llama3.2:latest"
$ws.Range("C6").Value = "The maximum number of characters in a single text entry is 32000."
